$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 3 new columns (Authors, Journal, Year) before old column C ---
$ws.Columns("C:E").Insert()

# --- Step 2: insert 1 new column (References (Top 5)) before the Stata Code column ---
$ws.Columns("R:R").Insert()

# --- Step 3: (re)write every cell with final content ---
# Row 1
$ws.Range("A1").Value = 'Filename'
$ws.Range("B1").Value = 'Title'
$ws.Range("C1").Value = 'Authors'
$ws.Range("D1").Value = 'Journal'
$ws.Range("E1").Value = 'Year'
$ws.Range("F1").Value = 'Background'
$ws.Range("G1").Value = 'Significance'
$ws.Range("H1").Value = 'Logic'
$ws.Range("I1").Value = 'Methodology'
$ws.Range("J1").Value = 'Conclusions'
$ws.Range("K1").Value = 'Dep. Var'
$ws.Range("L1").Value = 'Indep. Var'
$ws.Range("M1").Value = 'Mechanism'
$ws.Range("N1").Value = 'Instrumental'
$ws.Range("O1").Value = 'Controls'
$ws.Range("P1").Value = 'Data Source'
$ws.Range("Q1").Value = 'Measurements'
$ws.Range("R1").Value = 'References (Top 5)'
$ws.Range("S1").Value = 'Stata Code'

# Row 2
$ws.Range("A2").Value = '1-QJE-原神论文.pdf'
$ws.Range("B2").Value = '数字分心与同伴影响：移动应用使用对学术和劳动市场结果的影响'
$ws.Range("C2").Value = 'Panle Jia Barwick, Siyu Chen, Chao Fu, Teng Li'
$ws.Range("D2").Value = 'Quarterly Journal of Economics'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2024'
$ws.Range("F2").Value = '随着对手机过度使用的担忧日益增加，尤其是在青少年和年轻人中，本研究提供了首个关于个人及同伴移动应用使用对学术表现、身体健康和劳动市场结果影响的综合证据。研究利用中国一所大学的行政数据和手机记录，通过随机室友分配和政策冲击，探讨了个人及同伴的移动应用使用对学术成绩、身体健康和劳动市场结果的影响。'
$ws.Range("G2").Value = '本研究不仅在理论上增进了对数字分心和同伴效应的理解，而且在实践上为制定相关政策提供了依据，如限制青少年手机使用时间的政策可能会对他们未来的劳动市场结果产生积极影响。'
$ws.Range("H2").Value = '研究首先通过随机室友分配来估计同伴效应，然后利用政策变化作为工具变量来分离行为溢出效应和情境同伴效应，最后探讨了移动应用使用对学术成绩和劳动市场结果的影响。'
$ws.Range("I2").Value = '研究者首先利用随机室友分配来估计同伴效应，然后通过政策变化作为工具变量来分离行为溢出效应和情境同伴效应。接着，研究者分析了个人及同伴的移动应用使用对学术成绩和劳动市场结果的影响，并探讨了这些影响的异质性。最后，研究者利用高频GPS数据和在线调查数据来探讨这些影响背后的机制。'
$ws.Range("J2").Value = '[''移动应用使用具有传染性，室友的移动应用使用增加会导致个人使用增加。'', ''高移动应用使用对所有测量结果都是有害的，包括GPA和工资。'', ''室友的移动应用使用对个人的GPA和工资有直接和间接效应。'', ''限制未成年人游戏时间的政策可以提高大学生的初始工资。'', ''高频GPS数据显示，高应用使用挤出了学习时间和增加了迟到及缺课。'']'
$ws.Range("K2").Value = '学术成绩（GPA）、劳动市场结果（工资）'
$ws.Range("L2").Value = '个人及室友的移动应用使用时间'
$ws.Range("M2").Value = '时间分配（学习时间、迟到和缺课）'
$ws.Range("N2").Value = '未成年人游戏限制政策、爆款游戏“原神”发布时间'
$ws.Range("O2").Value = '年龄、农村居住、高中理科/文科轨迹、大学入学考试成绩、房价（家庭财富代理变量）'
$ws.Range("P2").Value = '中国一所大学的行政记录和手机使用数据，以及通过GPS系统和在线调查收集的补充数据。'
$ws.Range("Q2").Value = '通过手机记录测量移动应用使用时间，通过行政记录测量学术成绩和劳动市场结果。'
$ws.Range("R2").Value = 'Abdulkadiroğlu, Atila, Joshua Angrist, and Parag Pathak, “The elite illusion: Achievement effects at Boston and New York exam schools,” Econometrica, 2014.
Allcott, Hunt, Luca Braghieri, Sarah Eichmeyer, and Matthew Gentzkow, “The welfare effects of social media,” American Economic Review, 2020.
Brock, William A and Steven N Durlauf, “A multinomial choice model with social interactions,” in “The Economy as an Evolving Complex System III,” Oxford University Press, New York, 2006.
Sacerdote, Bruce, “Peer effects with random assignment: Results for Dartmouth roommates,” The Quarterly Journal of Economics, 2001.'
$ws.Range("S2").Value = '* 以二元游戏限制政策作为工具变量的2SLS回归
ivregress 2sls y x (z = iv), first
* 以爆款游戏发布作为工具变量的2SLS回归
ivregress 2sls y x (z = iv), first'

# Row 3
$ws.Range("A3").Value = '2-中国造船业的产业政策实施.pdf'
$ws.Range("B3").Value = '工业政策实施：来自中国造船业的实证证据'
$ws.Range("C3").Value = 'Panle Jia Barwick, Myrto Kalouptsidi, Nahim Bin Zahur'
$ws.Range("D3").Value = '未提供'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2023'
$ws.Range("F3").Value = '工业政策在全球范围内被广泛使用。设计和实施这些政策是一项复杂的任务。本文评估了不同工业政策工具的长期绩效，包括生产补贴、投资补贴、进入补贴和整合政策。研究以中国造船业为目标，该行业在21世纪初占全球产量的不到10%，但在国家五年计划期间获得多项政策干预后，迅速成为全球领先的船舶生产国。然而，这一增长伴随着大量新企业的涌入，加剧了行业分散和低产能利用率问题。'
$ws.Range("G3").Value = '本文的理论意义在于提供了工业政策设计的微观基础，强调了企业异质性、商业周期和企业成本结构在政策设计中的关键作用。实践意义在于，通过评估中国造船业的工业政策，为其他国家制定和实施工业政策提供了经验和教训。'
$ws.Range("H3").Value = '研究首先构建了一个包含动态和企业异质性的结构模型，然后使用1998至2014年的企业层面数据进行估计。研究的主要挑战是缺乏关于政府补贴性质的信息，因此通过比较政策实施前后中国企业的成本结构来恢复补贴的规模。接着，通过反事实分析评估中国工业政策的长期影响，并探讨不同政策工具的相对绩效。'
$ws.Range("I3").Value = '研究使用动态模型估计企业进入、退出和资本投资的决策。模型中，企业在每个时期进行库诺竞争，选择生产量，然后决定是否退出以及在继续运营的情况下进行多少投资。潜在进入者基于预期的终身利润流和进入成本做出一次性进入决策。通过估计模型，研究评估了中国工业政策对全球造船业的长期影响，并探讨了不同政策工具的相对绩效。'
$ws.Range("J3").Value = '[''中国政府的补贴政策大幅提高了中国在全球造船业的市场份额，但大部分增长来自于从其他国家企业抢夺业务。'', ''不同政策工具的绩效差异显著，生产和投资补贴相对于进入补贴更为有效。'', ''工业政策的有效性受到商业周期和企业效率异质性的影响，反周期政策的表现优于顺周期政策。'', ''政策通过整合和针对低成本企业进行补贴可以显著减少扭曲，但中国政府的实际白名单政策偏向于国有企业，牺牲了最有效率的企业。'']'
$ws.Range("K3").Value = '企业生产量、企业进入和退出决策、企业投资'
$ws.Range("L3").Value = '生产补贴、投资补贴、进入补贴、整合政策'
$ws.Range("M3").Value = '市场需求、企业成本结构、企业异质性'
$ws.Range("N3").Value = '无'
$ws.Range("O3").Value = '企业资本存量、企业积压订单、企业年龄、企业所有权状态、地区、企业规模'
$ws.Range("P3").Value = '数据来源包括Clarksons的全球船厂季度信息数据库和中国国家统计局的制造业企业年度数据库。'
$ws.Range("Q3").Value = '使用企业层面数据，包括产量、资本存量、特征以及船舶市场价格。'
$ws.Range("R3").Value = 'Johnson, 1982
Ito, 1992
Amsden, 1989
Lane, 2022
Stiglitz and Lin, 2013'
$ws.Range("S3").Value = '* 由于Stata代码较长，以下为简化示例
 reg y x1 x2, vce(cluster clustervar)
 ivreg2 y (x1 x2) z (x3 x4), first
 xtset panelvar timevar
 xtreg y x1 x2, fe
 * 更多复杂的动态面板数据模型和结构估计代码将根据具体方法论详细编写'

# Row 4
$ws.Range("A4").Value = '3-中国企业注册改革.pdf'
$ws.Range("B4").Value = '企业注册改革、市场周转与效率：中国的企业注册改革'
$ws.Range("C4").Value = 'Panle Jia Barwick, Luming Chen, Shanjun Li, Xiaobo Zhang'
$ws.Range("D4").Value = '八月 2024'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2024'
$ws.Range("F4").Value = '尽管各国普遍存在市场准入管制，但对这些管制如何影响企业动态和生产力的全面评估却相对缺乏。本文考察了2012-2014年中国广东省的试点项目（后来成为全国政策），旨在降低企业注册成本，鼓励创业活动。利用企业注册和年度报告的行政数据，分析显示改革使制造业部门的企业进入增加了25%，退出增加了8.7%。改革后进入者的生产力比改革前提高了1.1%，可能原因是金融约束的放松和竞争的加剧。'
$ws.Range("G4").Value = '本文不仅填补了有关进入管制对企业动态和生产力影响的文献空白，而且提供了中国大规模政策变化对企业层面影响的证据。研究结果对于理解市场准入改革的经济效应具有重要的理论和实践意义。'
$ws.Range("H4").Value = '文章首先介绍了研究背景和数据来源，然后构建了实证模型，接着讨论了企业进入、退出、规模和生产力的结果，并总结了稳健性检验。最后，探讨了影响企业生产力提升的潜在渠道和宏观影响。'
$ws.Range("I4").Value = '本文采用差分法(DID)策略，利用中国广东省企业注册改革的渐进式推广，通过比较改革前后企业进入和退出的变化来评估改革的影响。同时，通过三重差分设计来评估不同行业和企业类型对改革的异质性反应。'
$ws.Range("J4").Value = '[''企业注册改革显著增加了企业进入率和退出率，提高了市场周转率。'', ''改革后新进入企业的生产力提高了1.1%，表明更激烈的市场竞争和更优质的创业者进入市场。'', ''改革对不同类型企业和行业的效应存在异质性，私企和放松管制程度较高的行业受益更多。'', ''改革对提高企业生产力和宏观经济增长具有积极影响。'']'
$ws.Range("K4").Value = '因变量包括企业进入数量、退出情况、企业规模（如营收和资本）和企业生产力。'
$ws.Range("L4").Value = '自变量是改革实施的二元变量，表示是否在改革后的时间点。'
$ws.Range("M4").Value = '机制变量包括企业融资来源、企业家特征和企业股权结构。'
$ws.Range("N4").Value = '工具变量是改革实施的时间，利用城市间改革实施的时间差异作为工具变量。'
$ws.Range("O4").Value = '控制变量包括企业年龄、企业所有权类型、城市-行业固定效应、城市-年份固定效应和行业-年份固定效应等。'
$ws.Range("P4").Value = '数据来源包括广东省2009年1月至2016年12月的企业注册数据库、企业年度报告数据库和2018年中国创新与创业企业调查(ESIEC)。'
$ws.Range("Q4").Value = '企业进入和退出通过注册记录和年度报告数据进行测量。企业规模通过营收和资本来衡量，企业生产力通过Aw et al. (2011)和Peters et al. (2017)的方法进行结构性估计。'
$ws.Range("R4").Value = 'Aghion, P., Burgess, R., Redding, S. J., and Zilibotti, F. (2008). The unequal effects of liberalization: Evidence from dismantling the license raj in india. American Economic Review, 98(4):1397–1412.
Alfaro, L. and Chari, A. (2014). Deregulation, misallocation, and size: Evidence from india. The Journal of Law and Economics, 57(4):897–936.
Asturias, J., Hur, S., Kehoe, T. J., and Ruhl, K. J. (2023). Firm entry and exit and aggregate growth. American Economic Journal: Macroeconomics, 15(1):48–105.
Aw, B. Y., Roberts, M. J., and Xu, D. Y. (2011). R&D investment, exporting, and productivity dynamics. American Economic Review, 101(4):1312–44.
Baily, M. N., Hulten, C., Campbell, D., Bresnahan, T., and Caves, R. E. (1992). Productivity dynamics in manufacturing plants. Brookings Papers on Economic Activity. Microeconomics, 1992:187–267.'
$ws.Range("S4").Value = ""

# Row 5
$ws.Range("A5").Value = '4-城市交通政策的经济影响.pdf'
$ws.Range("B5").Value = '城市交通政策的经济影响与均衡排序'
$ws.Range("C5").Value = 'Panle Jia Barwick, Shanjun Li, Andrew Waxman, Jing Wu, Tianli Xia'
$ws.Range("D5").Value = 'Journal of Urban Economics'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2024'
$ws.Range("F5").Value = '交通在塑造城市空间结构和经济活动组织中起着关键作用。许多发展中国家由于快速城市化和机动化，加之基础设施不足，造成了前所未有的交通拥堵和严重的经济后果。全球各地政府实施了一系列政策，包括限行、公共交通投资、拥堵收费和汽油税等，以应对这些挑战。这些政策短期内在缓解拥堵方面的有效性，关键在于出行方式的可替代性和出行需求对通勤成本变化的敏感性。从中期到长期来看，这些政策可能会通过家庭调整居住地点对城市空间结构产生更广泛的影响。这种排序响应反过来可能影响交通政策在减少拥堵方面的有效性。此外，许多解决拥堵的政策具有分配后果。例如，收取通行费可能会加剧公平问题，因为低收入家庭在交通上的支出占收入比例更大。本文旨在理解城市交通政策的效率和公平影响，同时考虑排序响应和内生拥堵。'
$ws.Range("G5").Value = '本文的理论意义在于发展和估计了一个包含偏好异质性的住宅排序均衡模型，允许通过内生拥堵在住房地点和通勤决策之间进行一般均衡反馈。实践意义在于，模型能够预测在不同交通政策下的新均衡结果，包括出行方式选择、家庭位置、拥堵水平、房价和社会福利。'
$ws.Range("H5").Value = '研究逻辑首先通过家庭旅行调查估计出行时间和货币成本的异质性偏好，然后利用估计参数和工作地点构建每个通勤者的“通勤便利”属性，并在第二步中估计家庭对房产属性的偏好。研究考虑了观察到的家庭异质性和解决潜在内生性问题。'
$ws.Range("I5").Value = '首先，使用家庭旅行调查数据估计出行模式选择的参数，包括时间和货币成本偏好。其次，利用估计的参数和工作地点构建“通勤便利”属性，并在住房需求模型中使用该属性估计家庭对房产属性的偏好。最后，通过模拟不同交通政策下的均衡住宅排序和交通结果，分解福利效应。'
$ws.Range("J5").Value = '[''不同政策对居住地点的空间模式和均衡房价的影响不同，甚至有时相反。'', ''住宅排序可能加强或削弱交通政策的拥堵减少潜力。'', ''交通政策在总体和不同收入群体中产生不同的福利含义。'']'
$ws.Range("K5").Value = '通勤模式选择、住房位置选择'
$ws.Range("L5").Value = '交通政策（限行、拥堵收费、地铁扩张）'
$ws.Range("M5").Value = '偏好异质性、内生拥堵'
$ws.Range("N5").Value = '车牌摇号中奖概率（用于工具变量）'
$ws.Range("O5").Value = '房产大小、年龄、街区地址、交易价格、日期、家庭收入、年龄、性别、婚姻状况、户籍状态等'
$ws.Range("P5").Value = '北京家庭旅行调查（BHTS）2010年和2014年的数据，以及2006-2014年北京居民的主要政府抵押贷款计划中的住房交易数据。'
$ws.Range("Q5").Value = '使用家庭旅行调查数据和住房交易数据，构建包括家庭和工作地点、出行模式选择、住房属性等的综合数据集。'
$ws.Range("R5").Value = 'Akbar, Prottoy, “Who Benefits from Faster Public Transit?”, 2020. Working Paper.
Akbar, Prottoy A, Victor Couture, Gilles Duranton, and Adam Storeygard, “Mobility and Congestion in Urban India,” American Economic Review, 2023, 113 (4), 1083–1111.
Anderson, Michael L, “Subways, strikes, and slowdowns: The impacts of public transit on traffic congestion,” The American Economic Review, 2014, 104 (9), 2763–2796.
Gu, Yizhen, Naijia Guo, Jing Wu, and Ben Zou, “Home Location Choices and the Gender Commute Gap,” Journal of Human Resources, 2021, pp. 1020–11263R2.
Hainmueller, Jens, “Entropy Balancing for Causal Effects: A Multivariate Reweighting Method to Produce Balanced Samples in Observational Studies,” Political Analysis, 2012, 20 (1), 25?46.'
$ws.Range("S5").Value = ""

# Row 6
$ws.Range("A6").Value = '5-污染信息的价值.pdf'
$ws.Range("B6").Value = '从雾到霾：污染信息的价值'
$ws.Range("C6").Value = 'Panle Jia Barwick, Shanjun Li, Liguo Lin, Eric Zou'
$ws.Range("D6").Value = 'Journal of Political Economy'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2022'
$ws.Range("F6").Value = '2013年，中国启动了一个具有里程碑意义的空气质量监测和实时数据公开项目，显著增加了公众获取和对污染信息的意识。该项目触发了一系列行为改变，如更强的避免户外污染暴露和增加防护产品支出。这些行为响应减轻了空气污染的死亡影响。保守估计表明，该计划的健康益处比成本高出一个数量级。研究突出了在经常遭受严重空气污染但缺乏污染数据收集和传播的发展中国家提高公众获取污染信息的好处。'
$ws.Range("G6").Value = '本研究在理论和实践上都具有重要意义。理论上，它填补了关于发展中国家公民日常生活中收集和传播污染信息影响的知识空白。实践上，它为发展中国家公共资金在改善信息基础设施和满足基本医疗、营养、教育需求之间的竞争提供了见解。'
$ws.Range("H6").Value = '研究逻辑分为两步。首先，记录公共空气质量信息广泛可用后污染信息可用性和公民意识的变化。其次，研究监测计划如何影响行为和健康结果，重点关注文献中的两个核心对象：避免行为和空气污染的健康影响。通过检查户外活动和健康结果如何响应短期污染波动来衡量这些对象，并估计监测计划如何改变这些波动。'
$ws.Range("I6").Value = '本研究采用差分法，利用中国空气质量监测计划在不同城市的分阶段实施作为自然实验。研究构建了一个综合数据库，包括社会对空气污染的认知、卫星空气质量测量、经济活动和健康结果等丰富结果，覆盖监测计划前后的时期。'
$ws.Range("J6").Value = '[''监测计划显著提高了公众获取污染信息的途径，显著增加了家庭对污染问题的意识。'', ''信息获取和公众意识的变化触发了一系列行为变化，包括增加污染避免和防御性购买。'', ''监测计划有可能减轻严重空气污染的一些毁灭性健康后果。'', ''监测计划在提高健康结果方面的成本效益比非常高，是发展中国家记录的最成功的环境政策之一。'']'
$ws.Range("K6").Value = '户外活动频率、健康结果、空气净化器销售量'
$ws.Range("L6").Value = '空气污染水平（以AOD表示）'
$ws.Range("M6").Value = '污染信息的可用性和公民意识'
$ws.Range("N6").Value = '无'
$ws.Range("O6").Value = '城市固定效应、年份固定效应、周年份固定效应、城市特定时间趋势'
$ws.Range("P6").Value = '包括社会对空气污染的认知、卫星空气质量测量、经济活动和健康结果等丰富结果，覆盖监测计划前后的时期。'
$ws.Range("Q6").Value = '使用AOD作为空气污染的代理，通过MODIS算法从NASA卫星获取。'
$ws.Range("R6").Value = 'Stigler, 1961
Ashenfelter and Greenstone, 2004
Murphy and Topel, 2006
Barwick et al., 2020
Greenstone et al., 2022'
$ws.Range("S6").Value = '* 差分法估计
reg Outcomect Pollutionct##Postct, vce(cluster city)
* 事件研究法估计
reg Outcomect Pollutionct##c.q(), vce(cluster city)'

Write-Host "edit applied"
